$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Company"
$ws.Range("D1").Value = "Message"
$ws.Range("E1").Value = "Sent"

$message = "Happened to visit your company's website,excited to learn more about this colMail app."
$sent = "0:52 18-06-2023"

# Row 2 - new entry (Arpit)
$ws.Range("A2").Value = "verma.arpit078@gmail.com"
$ws.Range("B2").Value = "Arpit"
$ws.Range("C2").Value = "Google"
$ws.Range("D2").Value = $message
$ws.Range("E2").Value = $sent

# Row 3 - new entry (college)
$ws.Range("A3").Value = "2021eeb1156@iitrpr.ac.in"
$ws.Range("B3").Value = "college"
$ws.Range("C3").Value = "Google"
$ws.Range("D3").Value = $message
$ws.Range("E3").Value = $sent

# Row 4 - was previously row 2 (bla)
$ws.Range("A4").Value = "blabbla429@gmail.com"
$ws.Range("B4").Value = "bla"
$ws.Range("C4").Value = "Google"
$ws.Range("D4").Value = $message
$ws.Range("E4").Value = $sent

# Row 5 - was previously row 3 (server)
$ws.Range("A5").Value = "server1078@outlook.com"
$ws.Range("B5").Value = "server"
$ws.Range("C5").Value = "Google"
$ws.Range("D5").Value = $message
$ws.Range("E5").Value = $sent
